# edit.ps1
# Applies:
#   1. Insert "This role had some sub-projects, so provide a summary of the
#      role here." after the Summary heading that precedes "Some additional
#      bullets" (Another Role Name/Title section).
#   2. Remove the <w:ind w:left="720"/> direct formatting from the
#      "List some things that were done internally, if it applies" bullet.
#   3. Insert "This area is optional." after the Summary heading that
#      precedes "Technical Skills" (Generic Job Title section).
#   4. Insert "This area is optional." after the Summary heading that
#      precedes "Responsibilities" (the final role / Another Generic
#      Company Name section).

$d = $word.ActiveDocument

$wordNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Get-ParaText($index) {
    # Paragraph Range.Text includes a trailing paragraph-mark character that
    # does not print visibly, so trim it for reliable comparisons.
    return $d.Paragraphs.Item($index).Range.Text.TrimEnd()
}

function Find-ParaByNext($targetText, $nextContains) {
    for ($i = 1; $i -lt $d.Paragraphs.Count; $i++) {
        $t = Get-ParaText $i
        if ($t -eq $targetText) {
            $nxt = Get-ParaText ($i + 1)
            if ($nxt.Contains($nextContains)) {
                return $i
            }
        }
    }
    return -1
}

function Find-ParaByContains($substring) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = Get-ParaText $i
        if ($t.Contains($substring)) {
            return $i
        }
    }
    return -1
}

function Insert-ParagraphAfterHeading4Summary($nextContains, $newParaText) {
    $idx = Find-ParaByNext "Summary" $nextContains
    if ($idx -eq -1) {
        throw "Could not locate Summary heading before '$nextContains'"
    }
    $p = $d.Paragraphs.Item($idx)
    $r = $p.Range
    $xml = "<w:p $wordNs><w:pPr><w:pStyle w:val='Heading4'/></w:pPr><w:r><w:t>Summary</w:t></w:r></w:p><w:p $wordNs><w:r><w:t>$newParaText</w:t></w:r></w:p>"
    $r.InsertXML($xml) | Out-Null
}

# --- Change 1 ---------------------------------------------------------
Insert-ParagraphAfterHeading4Summary "Some additional bullets" "This role had some sub-projects, so provide a summary of the role here."

# --- Change 2 -----------------------------------------------------------
$idx2 = Find-ParaByContains "List some things that were done internally"
if ($idx2 -eq -1) {
    throw "Could not locate the 'List some things...' bullet paragraph"
}
$p2 = $d.Paragraphs.Item($idx2)
$r2 = $p2.Range
$xml2 = "<w:p $wordNs><w:pPr><w:pStyle w:val='ListBullet'/></w:pPr><w:r><w:t>List some things that were done internally, if it applies</w:t></w:r></w:p>"
$r2.InsertXML($xml2) | Out-Null

# --- Change 3 -------------------------------------------------------------
Insert-ParagraphAfterHeading4Summary "Technical Skills" "This area is optional."

# --- Change 4 -------------------------------------------------------------
Insert-ParagraphAfterHeading4Summary "Responsibilities" "This area is optional."

Write-Host "All edits applied."
